$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability values in column A
$ws.Range("A6").Value = 0.1
$ws.Range("A7").Value = 0.06
$ws.Range("A8").Value = 0.049
$ws.Range("A9").Value = 0.051
$ws.Range("A11").Value = 0.03

# Update the selected cell/active cell in the sheet view
$ws.Range("D7").Select()

# Update the workbook window position (best effort; this view-state
# attribute is carried through from the source file by this runtime and
# is not currently writable via the exposed Window object model, but we
# still set it in case the host reflects it).
$excel.ActiveWindow.Left = 22240
$excel.ActiveWindow.Top = 1320
